# Generate Report for Handback
# Regenerates the handback-status report with a fresh pair of UUID-named
# source files and fresh handoff/handback timestamps. Both zh-cn and de-de
# now resolve to the same hashed xlf name for the second row as for the
# first (the underlying content hash collided between the two files in this
# run), so several columns that used to differ between rows 2 and 3 now
# share the same value/cell text.

$wb = $excel.ActiveWorkbook

# ---- new literal values -----------------------------------------------
$mdFile1  = "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.md"
$mdFile2  = "ffffb7c21414-fd37-4e19-8bff-f3903dd8ea13.md"
$xlfZhCn  = "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.zh-cn.xlf"
$xlfDeDe  = "06b6d0d2-93fa-49a2-ad9b-85d379c57fe0.776ad46e8b0ea7fc280c238c0a86fdfc835cfa59.de-de.xlf"

$zhCnOffTime  = "2016-03-12 12:47:41"
$zhCnBackTime = "2016-03-12 12:47:58"
$deDeOffTime  = "2016-03-12 12:47:44"
$deDeBackTime = "2016-03-12 12:48:03"

# ---- Sheet 1: Overview --------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

$wsOverview.Range("A2").Value2 = $mdFile1
$wsOverview.Range("A3").Value2 = $mdFile2

$i = 0
foreach ($h in $wsOverview.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) {
        $h.TextToDisplay = $mdFile1
    } elseif ($i -eq 2) {
        $h.TextToDisplay = $mdFile2
    }
}

# ---- Sheet 2: zh-cn -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Range("A2").Value2 = $mdFile1
$wsZhCn.Range("D2").Value2 = $xlfZhCn
$wsZhCn.Range("E2").Value2 = $zhCnOffTime
$wsZhCn.Range("F2").Value2 = $mdFile1
$wsZhCn.Range("G2").Value2 = $xlfZhCn
$wsZhCn.Range("H2").Value2 = $zhCnBackTime

$wsZhCn.Range("A3").Value2 = $mdFile2
$wsZhCn.Range("D3").Value2 = $xlfZhCn
$wsZhCn.Range("E3").Value2 = $zhCnOffTime
$wsZhCn.Range("F3").Value2 = $mdFile2
$wsZhCn.Range("G3").Value2 = $xlfZhCn
$wsZhCn.Range("H3").Value2 = $zhCnBackTime

# Hyperlinks appear, in xlsx order, as: A2, B2, D2, F2, G2, A3, B3, D3, F3, G3
$i = 0
foreach ($h in $wsZhCn.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) {
        $h.TextToDisplay = $mdFile1
    } elseif ($i -eq 3) {
        $h.TextToDisplay = $xlfZhCn
    } elseif ($i -eq 4) {
        $h.TextToDisplay = $mdFile1
    } elseif ($i -eq 5) {
        $h.TextToDisplay = $xlfZhCn
    } elseif ($i -eq 6) {
        $h.TextToDisplay = $mdFile2
    } elseif ($i -eq 8) {
        $h.TextToDisplay = $xlfZhCn
    } elseif ($i -eq 9) {
        $h.TextToDisplay = $mdFile2
    } elseif ($i -eq 10) {
        $h.TextToDisplay = $xlfZhCn
    }
}

# ---- Sheet 3: de-de -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Range("A2").Value2 = $mdFile1
$wsDeDe.Range("D2").Value2 = $xlfDeDe
$wsDeDe.Range("E2").Value2 = $deDeOffTime
$wsDeDe.Range("F2").Value2 = $mdFile1
$wsDeDe.Range("G2").Value2 = $xlfDeDe
$wsDeDe.Range("H2").Value2 = $deDeBackTime

$wsDeDe.Range("A3").Value2 = $mdFile2
$wsDeDe.Range("D3").Value2 = $xlfDeDe
$wsDeDe.Range("E3").Value2 = $deDeOffTime
$wsDeDe.Range("F3").Value2 = $mdFile2
$wsDeDe.Range("G3").Value2 = $xlfDeDe
$wsDeDe.Range("H3").Value2 = $deDeBackTime

# Hyperlinks appear, in xlsx order, as: A2, B2, D2, F2, G2, A3, B3, D3, F3, G3
$i = 0
foreach ($h in $wsDeDe.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) {
        $h.TextToDisplay = $mdFile1
    } elseif ($i -eq 3) {
        $h.TextToDisplay = $xlfDeDe
    } elseif ($i -eq 4) {
        $h.TextToDisplay = $mdFile1
    } elseif ($i -eq 5) {
        $h.TextToDisplay = $xlfDeDe
    } elseif ($i -eq 6) {
        $h.TextToDisplay = $mdFile2
    } elseif ($i -eq 8) {
        $h.TextToDisplay = $xlfDeDe
    } elseif ($i -eq 9) {
        $h.TextToDisplay = $mdFile2
    } elseif ($i -eq 10) {
        $h.TextToDisplay = $xlfDeDe
    }
}

Write-Output "Handback report regenerated"
